# Automate creation of statewide SNOTEL layer, issue #147
#
# Adds a "Poudre RiverFest" entry to the Data sheet of the education
# organizations inventory, inserted (alphabetically, by Organization name)
# just above "South Platte Basin".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Record the existing Website-column hyperlinks (row + target URL) before
# the sheet is touched. Inserting a row shifts the underlying cell values
# down but does not renumber the <hyperlink> ref's already on the sheet, so
# we rebuild them by hand afterwards, in their original relationship order.
$links = @()
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $h = $ws.Hyperlinks.Item($i)
    $row = $h.Range.Row
    $val = $h.Range.Value()
    $links += , @($row, $val)
}

# Hyperlinks.Delete() clears the sheet's whole hyperlink collection (it is
# not per-link), which is fine here since every one of them is being
# re-created below at its shifted location.
$ws.Hyperlinks.Delete()

# Insert a new row at 11 (just above "South Platte Basin"), pushing that
# row and everything below it down by one.
$ws.Rows.Item(11).Insert()

# Fill in the new "Poudre RiverFest" row. Column order matches how the
# entry was authored: name, website text, description, programs, type,
# InBasin flag, then the coordinates.
$ws.Cells.Item(11, 1).Value = "Poudre RiverFest"
$ws.Cells.Item(11, 5).Value = "https://poudreriverfest.org/"
$ws.Cells.Item(11, 5).Style = "Hyperlink"
$ws.Cells.Item(11, 3).Value = "Education. Restoration. Celebration."
$ws.Cells.Item(11, 4).Value = "Annual Poudre River festival with exhibits."
$ws.Cells.Item(11, 2).Value = "Nonprofit"
$ws.Cells.Item(11, 6).Value = "Yes"
$ws.Cells.Item(11, 7).Value = -105.06962
$ws.Cells.Item(11, 8).Value = 40.59264

# Re-create the original hyperlinks, shifting any that were at/after the
# insertion point down by one row to track their (moved) cell content.
foreach ($l in $links) {
    $row = $l[0]
    $url = $l[1]
    if ($row -ge 11) { $row = $row + 1 }
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $url) | Out-Null
}

# Hyperlinks.Add() stamps a freshly duplicated "Hyperlink" style onto each
# cell it touches instead of reusing the workbook's existing one; reapply
# the named style across the whole Website column so it stays on the
# original shared style rather than growing the style table.
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 5).Style = "Hyperlink"
}

[void]$ws.Range("I11").Select()

Write-Output "Added 'Poudre RiverFest' to the Data sheet."
